$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.511.29"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.845.30"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("D5").Value = "'263.45"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5213"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "'0.3236"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").Value = "'0.06792"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "'18.71"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "'0.7756"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'0.07770"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "1.859.26"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "'88.25"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").Value = "'5.017"
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'13.95"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").Value = "'0.000007969"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "26.546.47"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "2.087.35"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "'4.616"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").Value = "'9.445"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").Value = "'5.985"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "'142.92"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  -8.55%  "
$ws.Range("D27").Value = "'1.677"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").Value = "'17.00"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'111.51"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'4.168"
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "'0.08733"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'4.102"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "'0.04823"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").Value = "'0.7203"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").Value = "'2.858"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "'3.088"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'0.01795"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").Value = "'2.203"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").Value = "'0.4841"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").Value = "'111.09"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").Value = "'0.8866"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "'6.040"
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'7.614"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "'0.4181"
$ws.Range("E46").Value = "  -2.07%  "
$ws.Range("D47").Value = "'0.05889"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "'9.003"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "'0.1235"
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").Value = "'34.90"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "'0.8884"
$ws.Range("E51").Value = "  +4.19%  "
